$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old wide data range (A1:M3) and reset row heights before rewriting the shrunk table
$ws.Range("A1:M3").Clear()
$ws.Rows("1:3").EntireRow.AutoFit()

# Header row
$ws.Cells.Item(1, 1).Value = "Sample no."
$ws.Cells.Item(1, 2).Value = "Monitor 1"
$ws.Cells.Item(1, 3).Value = "Monitor 2"

# Units row
$ws.Cells.Item(2, 1).Value = "Units"
$ws.Cells.Item(2, 2).Value = "[micro-strain]"
$ws.Cells.Item(2, 3).Value = "[micro-strain]"

# Data row
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = -1234.0999999999999
$ws.Cells.Item(3, 3).Value = -1168

# Styling: header cells (B1:C1) right aligned
$ws.Range("B1:C1").HorizontalAlignment = -4152

# Styling: units cells (B2:C2) centered
$ws.Range("B2:C2").HorizontalAlignment = -4108

# Styling: data cells (B3:C3) number format with 2 decimals
$ws.Range("B3:C3").NumberFormat = "0.00"

# Selection / view state
$ws.Range("A1:C3").Select() | Out-Null

$wb.Save() | Out-Null
